$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so that numeric-looking
# strings (e.g. "1.001", "5.458") are not silently converted into floating
# point numbers by Excel's automatic type detection.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.853.90"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "1.908.33"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "312.86"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "0.5228"
$ws.Range("E7").Value = "  +6.05%  "

$ws.Range("D8").Value = "0.3792"
$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("D9").Value = "0.07240"
$ws.Range("E9").Value = "  -1.20%  "

$ws.Range("D10").Value = "21.34"
$ws.Range("E10").Value = "  +3.76%  "

$ws.Range("D11").Value = "0.9089"
$ws.Range("E11").Value = "  -0.59%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07643"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.901.56"
$ws.Range("E13").Value = "  +0.06%  "

$ws.Range("D14").Value = "5.458"
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").Value = "92.38"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "0.000008707"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").Value = "27.872.62"
$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").Value = "5.153"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").Value = "2.166.05"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("E23").Value = "  +1.05%  "

$ws.Range("D24").Value = "6.626"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").Value = "153.79"
$ws.Range("E25").Value = "  +0.12%  "

$ws.Range("D27").Value = "2.173"
$ws.Range("E27").Value = "  +1.22%  "

$ws.Range("D28").Value = "18.34"
$ws.Range("E28").Value = "  -0.27%  "

$ws.Range("D29").Value = "114.66"
$ws.Range("E29").Value = "  -0.86%  "

$ws.Range("D30").Value = "4.855"
$ws.Range("E30").Value = "  -0.92%  "

$ws.Range("D31").Value = "0.09016"
$ws.Range("E31").Value = "  +0.87%  "

$ws.Range("D32").Value = "4.880"
$ws.Range("E32").Value = "  +5.16%  "

$ws.Range("D33").Value = "3.176"
$ws.Range("E33").Value = "  -0.71%  "

$ws.Range("D34").Value = "1.232"
$ws.Range("E34").Value = "  +1.01%  "

$ws.Range("D35").Value = "0.7809"
$ws.Range("E35").Value = "  +1.87%  "

$ws.Range("D36").Value = "0.02094"
$ws.Range("E36").Value = "  +3.13%  "

$ws.Range("D37").Value = "2.623"
$ws.Range("E37").Value = "  +3.92%  "

$ws.Range("D38").Value = "3.080"
$ws.Range("E38").Value = "  +3.36%  "

$ws.Range("D39").Value = "1.093"
$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("D40").Value = "0.5560"
$ws.Range("E40").Value = "  +1.63%  "

$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").Value = "6.699"
$ws.Range("E42").Value = "  -3.10%  "

$ws.Range("D43").Value = "115.06"
$ws.Range("E43").Value = "  +2.25%  "

$ws.Range("D44").Value = "8.585"
$ws.Range("E44").Value = "  +0.63%  "

$ws.Range("D45").Value = "0.1514"

$ws.Range("D46").Value = "0.4812"

$ws.Range("D47").Value = "10.43"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("D48").Value = "0.9990"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("D49").Value = "1.620"
$ws.Range("E49").Value = "  -0.63%  "

$ws.Range("E50").Value = "  -0.76%  "

$ws.Range("D51").Value = "0.06000"
$ws.Range("E51").Value = "  -0.80%  "

# Restore the default (Normal) style on the Price column so that no stray
# cell-format metadata is introduced; the cells keep their text values.
$priceRange.Style = "Normal"
